# Kpi_General_Export.xlsx - "thêm đơn hàng trực tiếp vào kpi"
# Adds 4 new KPI metric rows (SKUDirectOrder, TotalDirectSalesAmount,
# TotalDirectQuantity, TotalDirectOrders) right below the existing
# outlined "Đơn hàng gián tiếp" (indirect order) block, mirroring the
# layout used for rows 7-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C..T (18 columns) and the per-metric suffixes used throughout
# the sheet (Name, 12 months, 4 quarters, Year).
$cols      = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
$suffixes  = @("Name","M1","M2","M3","M4","M5","M6","M7","M8","M9","M10","M11","M12","Q1","Q2","Q3","Q4","Y")
$metrics   = @("SKUDirectOrder","TotalDirectSalesAmount","TotalDirectQuantity","TotalDirectOrders")
$newRows   = @(14,15,16,17)

# First copy the row formatting from the nearest existing rows so the
# new rows reuse the same cell styles instead of creating fresh ones:
#  - row 13 (plain "left/vcenter" style) -> rows 14, 16, 17
#  - row 11 (the yearly "Y" style with #,##0 number format) -> row 15 (D:T)
#  - C15 keeps the plain style (like C13)
$ws.Range("C13:T13").Copy() | Out-Null
$ws.Range("C14:T14").PasteSpecial(-4122) | Out-Null
$ws.Range("C16:T16").PasteSpecial(-4122) | Out-Null
$ws.Range("C17:T17").PasteSpecial(-4122) | Out-Null

$ws.Range("C13").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null

$ws.Range("D11:T11").Copy() | Out-Null
$ws.Range("D15:T15").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Fill in the template placeholder text for each new metric row.
for ($m = 0; $m -lt $metrics.Length; $m++) {
    $metric = $metrics[$m]
    $row = $newRows[$m]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = $cols[$i] + $row
        $ws.Range($addr).Value = "{{KpiGenerals." + $metric + "." + $suffixes[$i] + "}}"
    }
}

# Match the author's final selection in the saved workbook.
$ws.Range("AB5").Select()
